$d = $word.ActiveDocument

# 1. "7+ years of IT industry experience" -> "8+ years ..."
#    Replace only the leading "7" (keep rest of the sentence/run layout as close as possible).
$d.Content.Find.Execute("7+ years of IT industry", $true, $false, $false, $false, $false, $true, 1, $false, "8+ years of IT industry", 2) | Out-Null

# 2. Tools/tech list: drop "Wire-Shark, " in front of GIT.
$d.Content.Find.Execute("Wire-Shark, GIT", $true, $false, $false, $false, $false, $true, 1, $false, "GIT", 2) | Out-Null

# 3. Tools/tech list: "Bamboo" -> "Jenkins"
$d.Content.Find.Execute(", Bamboo", $true, $false, $false, $false, $false, $true, 1, $false, ", Jenkins", 2) | Out-Null

# 4. "..., Kafka, Zookeeper." -> "..., Kafka, Zookeeper, Slf4j"
$d.Content.Find.Execute("Kafka, Zookeeper.", $true, $false, $false, $false, $false, $true, 1, $false, "Kafka, Zookeeper, Slf4j", 2) | Out-Null

# 5. Bullet rewrite: "Refactored code ..." -> "Implemented Lapse, Reinstatement and Maturity in D-PAS. "
$d.Content.Find.Execute("Refactored code by writing utility classes for pdf  generation and PRT communication.", $true, $false, $false, $false, $false, $true, 1, $false, "Implemented Lapse, Reinstatement and Maturity in D-PAS. ", 2) | Out-Null

# 6. Tech list: drop "Apache-Chain, "
$d.Content.Find.Execute("Java-SE, Apache-Chain, XML, Log4j", $true, $false, $false, $false, $false, $true, 1, $false, "Java-SE, XML, Log4j", 2) | Out-Null
